$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2325581395348837
$ws.Range("C2").Value = 0.5038759689922481
$ws.Range("J2").Value = 0.007751937984496124
$ws.Range("P2").Value = 0.1550387596899225
$ws.Range("S2").Value = 0.1007751937984496
$ws.Range("C3").Value = 0.01550387596899225
$ws.Range("J3").Value = 0.03875968992248062
$ws.Range("P3").Value = 0.7674418604651163
$ws.Range("S3").Value = 0.1782945736434109
$ws.Range("J4").Value = 0.08823529411764706
$ws.Range("P4").Value = 0.6470588235294118
$ws.Range("S4").Value = 0.2647058823529412
$ws.Range("P5").Value = 1
$ws.Range("B6").Value = 0.0380952380952381
$ws.Range("D6").Value = 0.01428571428571429
$ws.Range("F6").Value = 0.06190476190476191
$ws.Range("J6").Value = 0.2571428571428571
$ws.Range("O6").Value = 0.01904761904761905
$ws.Range("Q6").Value = 0.1619047619047619
$ws.Range("R6").Value = 0.05714285714285714
$ws.Range("S6").Value = 0.3904761904761905
$ws.Range("B7").Value = 0.08888888888888889
$ws.Range("D7").Value = 0.02777777777777778
$ws.Range("E7").Value = 0.005555555555555556
$ws.Range("F7").Value = 0.04444444444444445
$ws.Range("J7").Value = 0.1055555555555556
$ws.Range("O7").Value = 0.01111111111111111
$ws.Range("Q7").Value = 0.1611111111111111
$ws.Range("R7").Value = 0.03888888888888889
$ws.Range("S7").Value = 0.5166666666666667
$ws.Range("B8").Value = 0.06
$ws.Range("D8").Value = 0.008
$ws.Range("F8").Value = 0.046
$ws.Range("J8").Value = 0.11
$ws.Range("O8").Value = 0.02
$ws.Range("Q8").Value = 0.18
$ws.Range("R8").Value = 0.116
$ws.Range("S8").Value = 0.46
$ws.Range("B9").Value = 0.08230452674897119
$ws.Range("D9").Value = 0.01234567901234568
$ws.Range("E9").Value = 0.00411522633744856
$ws.Range("F9").Value = 0.06995884773662552
$ws.Range("J9").Value = 0.06995884773662552
$ws.Range("O9").Value = 0.01234567901234568
$ws.Range("Q9").Value = 0.1975308641975309
$ws.Range("R9").Value = 0.07407407407407407
$ws.Range("S9").Value = 0.4773662551440329
$ws.Range("B10").Value = 0.09597523219814241
$ws.Range("D10").Value = 0.01857585139318885
$ws.Range("E10").Value = 0.0007739938080495357
$ws.Range("F10").Value = 0.06888544891640867
$ws.Range("J10").Value = 0.1029411764705882
$ws.Range("O10").Value = 0.02012383900928793
$ws.Range("Q10").Value = 0.2105263157894737
$ws.Range("R10").Value = 0.0673374613003096
$ws.Range("S10").Value = 0.4148606811145511
$ws.Range("G11").Value = 0.157556270096463
$ws.Range("J11").Value = 0.1286173633440514
$ws.Range("K11").Value = 0.2315112540192926
$ws.Range("L11").Value = 0.4565916398713826
$ws.Range("S11").Value = 0.02572347266881029
$ws.Range("G12").Value = 0.695364238410596
$ws.Range("J12").Value = 0.2119205298013245
$ws.Range("K12").Value = 0.006622516556291391
$ws.Range("L12").Value = 0.03973509933774835
$ws.Range("S12").Value = 0.04635761589403974
$ws.Range("G13").Value = 0.7045454545454546
$ws.Range("J13").Value = 0.2045454545454546
$ws.Range("S13").Value = 0.09090909090909091
$ws.Range("F15").Value = 0.02512562814070352
$ws.Range("H15").Value = 0.1658291457286432
$ws.Range("I15").Value = 0.07035175879396985
$ws.Range("J15").Value = 0.3517587939698493
$ws.Range("K15").Value = 0.03015075376884422
$ws.Range("M15").Value = 0.01507537688442211
$ws.Range("O15").Value = 0.09045226130653267
$ws.Range("S15").Value = 0.2512562814070352
$ws.Range("F16").Value = 0.03821656050955414
$ws.Range("H16").Value = 0.1847133757961783
$ws.Range("I16").Value = 0.1146496815286624
$ws.Range("J16").Value = 0.4076433121019108
$ws.Range("K16").Value = 0.1082802547770701
$ws.Range("M16").Value = 0.01273885350318471
$ws.Range("O16").Value = 0.01273885350318471
$ws.Range("S16").Value = 0.1210191082802548
$ws.Range("F17").Value = 0.01279317697228145
$ws.Range("H17").Value = 0.1833688699360341
$ws.Range("I17").Value = 0.1023454157782516
$ws.Range("J17").Value = 0.417910447761194
$ws.Range("K17").Value = 0.08528784648187633
$ws.Range("M17").Value = 0.01066098081023454
$ws.Range("N17").Value = 0.002132196162046908
$ws.Range("O17").Value = 0.05543710021321962
$ws.Range("S17").Value = 0.1300639658848614
$ws.Range("F18").Value = 0.0273224043715847
$ws.Range("H18").Value = 0.2076502732240437
$ws.Range("I18").Value = 0.09289617486338798
$ws.Range("J18").Value = 0.366120218579235
$ws.Range("K18").Value = 0.09836065573770492
$ws.Range("M18").Value = 0.01639344262295082
$ws.Range("O18").Value = 0.06557377049180328
$ws.Range("S18").Value = 0.1256830601092896
$ws.Range("F19").Value = 0.008934707903780068
$ws.Range("H19").Value = 0.2219931271477663
$ws.Range("I19").Value = 0.1010309278350515
$ws.Range("J19").Value = 0.3738831615120275
$ws.Range("K19").Value = 0.1010309278350515
$ws.Range("M19").Value = 0.02199312714776632
$ws.Range("N19").Value = 0.0006872852233676976
$ws.Range("O19").Value = 0.04810996563573883
$ws.Range("S19").Value = 0.1223367697594502
